$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.000.46"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "3.386.78"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'573.11"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'137.07"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.384.21"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").Value = "'7.63"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D11").Value = "'0.123"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").Value = "3.964.05"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "'26.61"
$ws.Range("D16").Value = "3.386.49"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("D18").Value = "61.047.85"
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'5.87"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'13.88"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "'9.31"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "'375.24"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "'0.550"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").Value = "3.512.69"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("E28").Value = "  -6.98%  "
$ws.Range("E29").Value = "  +7.71%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'7.38"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").Value = "'8.06"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'23.40"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").Value = "'165.00"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'0.0767"
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("D41").Value = "'25.96"
$ws.Range("E41").Value = "  +4.89%  "
$ws.Range("D42").Value = "'1.75"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'41.93"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").Value = "'4.38"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("E47").Value = "  -4.29%  "
$ws.Range("D48").Value = "2.513.42"
$ws.Range("E48").Value = "  +7.22%  "
$ws.Range("D49").Value = "'23.66"
$ws.Range("E49").Value = "  +3.76%  "
$ws.Range("D50").Value = "'6.77"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").Value = "'2.40"
$ws.Range("E51").Value = "  +3.02%  "
